$wb = $excel.ActiveWorkbook

# Overview sheet
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E6").Value = "In Translation"
$wsOverview.Range("F6").Value = "In Translation"
$wsOverview.Range("G6").Value = "2016-10-20 08:37:43"

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C6").Value = "In Translation"
$wsZhCn.Range("H6").Value = "2016-10-20 08:37:32"

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C6").Value = "In Translation"
$wsDeDe.Range("H6").Value = "2016-10-20 08:37:43"
